# Update TPM-derived NATMI metrics for the Efna2-Epha5 LR-pairs sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ECs)
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.573649333333333
$ws.Range("H2").Value = 4.720948
$ws.Range("I2").Value = 0.162950296453897
$ws.Range("J2").Value = 0.1728167686459121
$ws.Range("M2").Value = 0.0237255
$ws.Range("N2").Value = 0.047451
$ws.Range("Q2").Value = 0.037335617258
$ws.Range("R2").Value = 0.224013703548
$ws.Range("S2").Value = 0.162950296453897
$ws.Range("T2").Value = 0.1728167686459121

# Row 3 (FAPs -> ECs)
$ws.Range("I3").Value = 0.4548971409363525
$ws.Range("J3").Value = 0.4824406931050072
$ws.Range("M3").Value = 0.0237255
$ws.Range("N3").Value = 0.047451
$ws.Range("Q3").Value = 0.104227276141
$ws.Range("R3").Value = 0.625363656846
$ws.Range("S3").Value = 0.4548971409363525
$ws.Range("T3").Value = 0.4824406931050072

# Row 4 (Inflammatory-Mac -> ECs)
$ws.Range("G4").Value = 1.149467
$ws.Range("H4").Value = 3.448401
$ws.Range("I4").Value = 0.11902651019285
$ws.Range("J4").Value = 0.1262334424813261
$ws.Range("M4").Value = 0.0237255
$ws.Range("N4").Value = 0.047451
$ws.Range("Q4").Value = 0.0272716793085
$ws.Range("R4").Value = 0.163630075851
$ws.Range("S4").Value = 0.11902651019285
$ws.Range("T4").Value = 0.1262334424813261

# Row 5 (MuSCs -> ECs)
$ws.Range("G5").Value = 1.6540555
$ws.Range("H5").Value = 3.308111
$ws.Range("I5").Value = 0.1712762992154535
$ws.Range("J5").Value = 0.1210979348516435
$ws.Range("M5").Value = 0.0237255
$ws.Range("N5").Value = 0.047451
$ws.Range("Q5").Value = 0.03924329376525
$ws.Range("R5").Value = 0.156973175061
$ws.Range("S5").Value = 0.1712762992154535
$ws.Range("T5").Value = 0.1210979348516435

# Row 6 (Resolving-Mac -> ECs)
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.8870146666666666
$ws.Range("H6").Value = 2.661044
$ws.Range("I6").Value = 0.09184975320144682
$ws.Range("J6").Value = 0.097411160916111
$ws.Range("M6").Value = 0.0237255
$ws.Range("N6").Value = 0.047451
$ws.Range("Q6").Value = 0.021044866474
$ws.Range("R6").Value = 0.126269198844
$ws.Range("S6").Value = 0.09184975320144682
$ws.Range("T6").Value = 0.097411160916111
